$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0224449634552002
$ws.Range("C2").Value = 0.04331727027893066
$ws.Range("D2").Value = 0.01298418045043945
$ws.Range("E2").Value = 0.02321395874023437
$ws.Range("F2").Value = 0.009221601486206054
$ws.Range("G2").Value = 0.1111621379852295
$ws.Range("H2").Value = 0.01633610725402832
$ws.Range("I2").Value = 0.03061103820800781
$ws.Range("J2").Value = 0.01359930038452148
$ws.Range("K2").Value = 0.02527627944946289
$ws.Range("L2").Value = 0.008940172195434571
$ws.Range("M2").Value = 0.02662134170532227
$ws.Range("B3").Value = 0.09064240455627441
$ws.Range("C3").Value = 0.02895145416259766
$ws.Range("D3").Value = 0.03159785270690918
$ws.Range("E3").Value = 0.02548084259033203
$ws.Range("F3").Value = 0.0145392894744873
$ws.Range("G3").Value = 0.01705532073974609
$ws.Range("H3").Value = 0.1523452281951904
$ws.Range("I3").Value = 0.04470338821411133
$ws.Range("J3").Value = 0.06529712677001953
$ws.Range("K3").Value = 0.02166919708251953
$ws.Range("L3").Value = 0.03138537406921386
$ws.Range("M3").Value = 0.01520314216613769
$ws.Range("B4").Value = 0.04771676063537598
$ws.Range("C4").Value = 0.03128142356872558
$ws.Range("D4").Value = 0.01005845069885254
$ws.Range("E4").Value = 0.01064667701721191
$ws.Range("F4").Value = 0.04741535186767578
$ws.Range("G4").Value = 0.01022868156433105
$ws.Range("H4").Value = 0.03005590438842774
$ws.Range("I4").Value = 0.02295050621032715
$ws.Range("J4").Value = 0.02708320617675781
$ws.Range("K4").Value = 0.02170171737670899
$ws.Range("L4").Value = 0.04835963249206543
$ws.Range("M4").Value = 0.009272575378417969
$ws.Range("B5").Value = 0.02702350616455078
$ws.Range("C5").Value = 0.0285923957824707
$ws.Range("D5").Value = 0.01271085739135742
$ws.Range("E5").Value = 0.01342740058898926
$ws.Range("H5").Value = 0.01565513610839844
$ws.Range("I5").Value = 0.02653164863586426
$ws.Range("J5").Value = 0.01306591033935547
$ws.Range("K5").Value = 0.0133671760559082
$ws.Range("B6").Value = 0.2495463371276855
$ws.Range("C6").Value = 0.09323277473449706
$ws.Range("D6").Value = 0.09794659614562988
$ws.Range("E6").Value = 0.04724011421203613
$ws.Range("F6").Value = 0.09229340553283691
$ws.Range("G6").Value = 0.04193358421325684
$ws.Range("H6").Value = 0.2245556354522705
$ws.Range("I6").Value = 0.0906646728515625
$ws.Range("J6").Value = 0.13083815574646
$ws.Range("K6").Value = 0.06593403816223145
$ws.Range("L6").Value = 0.0783452033996582
$ws.Range("M6").Value = 0.03139667510986328
